$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 used to hold the literal text "24/10/2022"; turn it into a real date
# serial (matching the date formatting already used by A1/A2) instead.
$ws.Range("A3").Value = 44858
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat

# B3 keeps pointing at the "24/10/2022" distribution file (hyperlink rId2
# already targets .../24_10_2022.xlsx?raw=true) - only the displayed text
# changes, from the old date label to the raw-file URL.
$ws.Range("B3").Value = "https://github.com/gandharvas/crs/blob/main/files/24_10_2022.xlsx?raw=true"

# New row for the 07/11/2022 cut-off date.
$ws.Range("A4").Value = 44872
$ws.Range("A4").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B4").Value = "https://github.com/gandharvas/crs/blob/main/files/07_11_2022.xlsx?raw=true"
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/gandharvas/crs/blob/main/files/07_11_2022.xlsx?raw=true") | Out-Null

# Hyperlinks.Add re-derives the cell's font/style; reapply the shared
# "Hyperlink" look used by B2/B3 so B4 matches its siblings exactly.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B5").Select()
